# Update countries & provincias Spain
#
# 1) "Somalia" moves earlier in the country list (its shared-string entry
#    is now located right before "Taiwan"), which shifts the country labels
#    for rows 112-115 down by one (Somalia, Taiwan, Reunion, Mayotte).
# 2) A handful of statistic rows (totals / active / recovered / critical /
#    deaths) are refreshed with newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country column (A) relabeling for rows 112-115 ------------------------
$ws.Range("A112").Value = "Somalia"
$ws.Range("A113").Value = "Taiwan"
$ws.Range("A114").Value = "Reunion"
$ws.Range("A115").Value = "Mayotte"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 964075
$ws.Range("C4").Value = 3424
$ws.Range("E4").Value = 791364
$ws.Range("G4").Value = 119
$ws.Range("H4").Value = 54375

# --- Row 16: Canada ---------------------------------------------------------
$ws.Range("B16").Value = 45791
$ws.Range("C16").Value = 437
$ws.Range("D16").Value = 16916
$ws.Range("E16").Value = 26386
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = 2489

# --- Row 33: Polonia ---------------------------------------------------------
$ws.Range("B33").Value = 11617
$ws.Range("C33").Value = 344
$ws.Range("E33").Value = 8817
$ws.Range("G33").Value = 11
$ws.Range("H33").Value = 535

# --- Row 34: Rumania ---------------------------------------------------------
$ws.Range("E34").Value = 7363
$ws.Range("G34").Value = 18
$ws.Range("H34").Value = 619

# --- Row 45: Chequia ---------------------------------------------------------
$ws.Range("B45").Value = 7387
$ws.Range("C45").Value = 35
$ws.Range("D45").Value = 2545
$ws.Range("E45").Value = 4622
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 220

# --- Row 83: Cuba -------------------------------------------------------------
$ws.Range("B83").Value = 1369
$ws.Range("C83").Value = 32
$ws.Range("D83").Value = 501
$ws.Range("E83").Value = 814
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 54

# --- Row 110: Jordania ---------------------------------------------------------
$ws.Range("B110").Value = 447
$ws.Range("C110").Value = 3
$ws.Range("D110").Value = 337
$ws.Range("E110").Value = 103

# --- Row 112: now Somalia (new entry) -------------------------------------
$ws.Range("B112").Value = 436
$ws.Range("C112").Value = 46
$ws.Range("D112").Value = 10
$ws.Range("E112").Value = 403
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 5
$ws.Range("H112").Value = 23

# --- Row 113: now Taiwan (shifted down, values unchanged) -----------------
$ws.Range("B113").Value = 429
$ws.Range("D113").Value = 281
$ws.Range("E113").Value = 142
$ws.Range("F113").Value = 0
$ws.Range("H113").Value = 6

# --- Row 114: now Reunion (shifted down, values unchanged) -----------------
$ws.Range("B114").Value = 417
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 300
$ws.Range("E114").Value = 117
$ws.Range("F114").Value = 2
$ws.Range("H114").Value = 0

# --- Row 115: now Mayotte (shifted down, values unchanged) -----------------
$ws.Range("B115").Value = 401
$ws.Range("C115").Value = 21
$ws.Range("D115").Value = 144
$ws.Range("E115").Value = 253
$ws.Range("F115").Value = 4
$ws.Range("H115").Value = 4

# --- Row 119: Mauricio ---------------------------------------------------------
$ws.Range("B119").Value = 332
$ws.Range("C119").Value = 1
$ws.Range("D119").Value = 299
$ws.Range("E119").Value = 24

# --- Row 159: San Martin (Parte Holandesa) ----------------------------------
$ws.Range("B159").Value = 74
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 33
$ws.Range("E159").Value = 28
$ws.Range("F159").Value = 7
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 13
